$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Objetivos:" (row 10) body text in Portuguese
$objetivosNew = 'Propiciar ao aluno um panorama geral da área de Física do Estado Sólido, com ênfase nas idéias fundamentais e conceitos gerais, como gás de elétron, excitações elementares, estrutura de bandas, etc. O curso deve ser rico em resultados experimentais que ilustrem princípios e comportamentos gerais dos sólidos (por exemplo, comportamento das grandezas físicas com a temperatura).'
$ws.Cells.Item(10, 2).Value = $objetivosNew
$ws.Cells.Item(10, 3).Value = $objetivosNew

# 2. Insert three new rows right after row 12 ("Docentes responsáveis:") to hold the
#    three professor names. They carry no label in column A, just the B/C text,
#    so clone the B/C formatting (style) from the still-untouched "Short syllabus:"
#    row (currently row 14, will become row 17 after the inserts) and then clear
#    out column A on the new rows.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(17, 2).Copy($ws.Cells.Item(13, 2))
$ws.Cells.Item(17, 3).Copy($ws.Cells.Item(13, 3))
$ws.Cells.Item(17, 2).Copy($ws.Cells.Item(14, 2))
$ws.Cells.Item(17, 3).Copy($ws.Cells.Item(14, 3))
$ws.Cells.Item(17, 2).Copy($ws.Cells.Item(15, 2))
$ws.Cells.Item(17, 3).Copy($ws.Cells.Item(15, 3))

$ws.Cells.Item(13, 1).Clear()
$ws.Cells.Item(14, 1).Clear()
$ws.Cells.Item(15, 1).Clear()

$ws.Cells.Item(13, 2).Value = '5840730 - Antonio Jefferson da Silva Machado'
$ws.Cells.Item(13, 3).Value = '5840730 - Antonio Jefferson da Silva Machado'
$ws.Cells.Item(14, 2).Value = '5840726 - Cristina Bormio Nunes'
$ws.Cells.Item(14, 3).Value = '5840726 - Cristina Bormio Nunes'
$ws.Cells.Item(15, 2).Value = '1341653 - Maria José Ramos Sandim'
$ws.Cells.Item(15, 3).Value = '1341653 - Maria José Ramos Sandim'

# 3. "Programa resumido:" is now row 16 - give it its real summary text
$programaResumidoNew = 'Estrutura e ligações cristalinas. Vibrações da rede, fônons e propriedades térmicas. Gás de Fermi de elétrons livres. Bandas de energia. Semicondutores. Metais e superfícies de Fermi.'
$ws.Cells.Item(16, 2).Value = $programaResumidoNew
$ws.Cells.Item(16, 3).Value = $programaResumidoNew

# 4. "Programa:" is now row 18 - give it its real full syllabus text
$programaNew = '¨ Estrutura dos cristais.¨ Difração em cristais e a rede recíproca.¨ Ligações em cristais: cristais iônicos e cristais covalentes¨ Constantes elásticas e ondas elásticas.¨ Vibrações de cristais. Fônons¨ Gás de Fermi: modelo do elétron livre; movimento em campos magnéticos.¨ Bandas de energia. Funções de Bloch.¨ Cristais semicondutores.'
$ws.Cells.Item(18, 2).Value = $programaNew
$ws.Cells.Item(18, 3).Value = $programaNew

# 5. "Método:" is now row 21 - real teaching method text
$metodoNew = 'Aulas expositivas, seminários e exercícios comentados.'
$ws.Cells.Item(21, 2).Value = $metodoNew
$ws.Cells.Item(21, 3).Value = $metodoNew

# 6. "Critério:" is now row 22 - real grading criteria text
$criterioNew = 'Média aritmética de duas provas com mesmo peso.'
$ws.Cells.Item(22, 2).Value = $criterioNew
$ws.Cells.Item(22, 3).Value = $criterioNew

# 7. "Norma de recuperação:" is now row 23 - real make-up exam rule text
$normaNew = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Cells.Item(23, 2).Value = $normaNew
$ws.Cells.Item(23, 3).Value = $normaNew

# 8. "Bibliografia:" is now row 24 - real bibliography text
$biblioNew = 'ASHCROFT, N. W. Solid State Physics. Saunders College. KITTEL, C. Introduction to Solid State Physics. John Wiley & Sons. BLAKEMORE, J. S. Solid State Physics, Cambridge University Press. WERT,  C. A.; THOMSON, R. B. Physics of Solids. McGraw-Hill Book Co. Ltda. 1968. ZIMAN, J. M. Principles of the theory of solids, Cambridge, 2nd ed., 1972. SUTTON, A. P.  Electronic Structure of Materials, Oxford Science Publications.'
$ws.Cells.Item(24, 2).Value = $biblioNew
$ws.Cells.Item(24, 3).Value = $biblioNew
